$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.021.26"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.53%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.914.34"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.82%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5033"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4027"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08241"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.107"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.04"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.75"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.922.74"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.425"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.282"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "92.11"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001097"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06508"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.31%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.944"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "30.043.88"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.201"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "22.36"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.134.58"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.67"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.277"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.91"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.135"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1038"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.014"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.794"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +4.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02442"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -1.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.345"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +3.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06424"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2161"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.95%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.828"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.201"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6464"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.41"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -4.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.221"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.07%  "
# Row 45: EnergySwap -> NEARProtocol (rank 43 now shows NEARProtocol)
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.198"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.53%  "

# Row 46: NEARProtocol -> EnergySwap (rank 44 now shows EnergySwap)
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.32"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5997"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.636"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.36"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.213"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "78.84"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.71%  "
